$d = $word.ActiveDocument

# 1. "Version 1." and "10" are two adjacent runs that together read
#    "Version 1.10". Re-typing the same text via Find/Replace collapses
#    them into a single run, matching the intended XML normalisation.
$d.Content.Find.Execute("Version 1.10", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version 1.10", 2) | Out-Null

# 2. Insert a new bulleted line after "Fixed preset failsafe bug..." that
#    documents the sBus baud-rate fix, re-using the same numbered list
#    (numId 6) as the paragraph above it.
$srcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Fixed preset failsafe bug.*") {
        $srcPara = $p
        break
    }
}

if ($srcPara -ne $null) {
    $srcTemplate = $srcPara.Range.ListFormat.ListTemplate

    $i = 0
    $targetIndex = -1
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text -like "Fixed preset failsafe bug.*") {
            $p.Range.InsertParagraphAfter()
            $targetIndex = $i + 1
            break
        }
    }

    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Style = "NoSpacing"
    $newPara.Range.Text = "sBus baud rate was slightly off."
    $newPara.Range.ListFormat.ApplyListTemplate($srcTemplate, $true)
}
